$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = $origStyle
}

$ws.Range('D2').Value = '66.894.59'
$ws.Range('E2').Value = '  +0.14%  '

$ws.Range('D3').Value = '3.116.84'
$ws.Range('E3').Value = '  +0.95%  '

$ws.Range('E4').Value = '  +0.03%  '

Set-TextValue 'D5' '577.14'
$ws.Range('E5').Value = '  -0.47%  '

Set-TextValue 'D6' '172.21'
$ws.Range('E6').Value = '  +2.53%  '

$ws.Range('E7').Value = '  +0.04%  '

$ws.Range('E8').Value = '  -0.42%  '

$ws.Range('E9').Value = '  -3.10%  '

$ws.Range('E10').Value = '  -1.66%  '

$ws.Range('E11').Value = '  -0.06%  '

$ws.Range('E12').Value = '  -1.36%  '

Set-TextValue 'D13' '37.11'
$ws.Range('E13').Value = '  +1.88%  '

$ws.Range('E14').Value = '  -1.21%  '

$ws.Range('D15').Value = '3.633.53'
$ws.Range('E15').Value = '  +0.93%  '

$ws.Range('D16').Value = '66.860.87'
$ws.Range('E16').Value = '  +0.09%  '

$ws.Range('E17').Value = '  -0.90%  '

$ws.Range('D18').Value = '3.116.89'
$ws.Range('E18').Value = '  +0.93%  '

Set-TextValue 'D19' '16.22'
$ws.Range('E19').Value = '  +0.72%  '

Set-TextValue 'D20' '474.95'
$ws.Range('E20').Value = '  +1.70%  '

Set-TextValue 'D21' '0.708'
$ws.Range('E21').Value = '  -0.74%  '

Set-TextValue 'D22' '7.90'
$ws.Range('E22').Value = '  +5.27%  '

Set-TextValue 'D23' '83.80'
$ws.Range('E23').Value = '  -0.05%  '

Set-TextValue 'D25' '2.27'
$ws.Range('E25').Value = '  -3.59%  '

$ws.Range('E26').Value = '  +0.46%  '

$ws.Range('E27').Value = '  +0.01%  '

Set-TextValue 'D28' '7.92'
$ws.Range('E28').Value = '  -0.74%  '

$ws.Range('E29').Value = '  -1.30%  '

Set-TextValue 'D30' '2.67'
$ws.Range('E30').Value = '  -0.03%  '

Set-TextValue 'D31' '28.50'
$ws.Range('E31').Value = '  +0.93%  '

$ws.Range('E32').Value = '  -0.51%  '

$ws.Range('D33').Value = '0.0₃0948'
$ws.Range('E33').Value = '  -7.14%  '

$ws.Range('E34').Value = '  -0.03%  '

Set-TextValue 'D35' '5.82'
$ws.Range('E35').Value = '  -1.18%  '

Set-TextValue 'D36' '0.974'
$ws.Range('E36').Value = '  -3.07%  '

Set-TextValue 'D37' '46.86'
$ws.Range('E37').Value = '  -0.39%  '

$ws.Range('B38').Value = 'OKB'
$ws.Range('C38').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue 'D38' '50.21'
$ws.Range('E38').Value = '  -0.14%  '

$ws.Range('B39').Value = 'Stacks'
$ws.Range('C39').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue 'D39' '2.06'
$ws.Range('E39').Value = '  -2.03%  '

Set-TextValue 'D40' '0.310'
$ws.Range('E40').Value = '  -2.59%  '

$ws.Range('E41').Value = '  +1.51%  '

Set-TextValue 'D42' '8.57'
$ws.Range('E42').Value = '  -1.28%  '

$ws.Range('D43').Value = '2.812.51'
$ws.Range('E43').Value = '  +1.11%  '

Set-TextValue 'D44' '382.29'
$ws.Range('E44').Value = '  -0.07%  '

$ws.Range('E45').Value = '  -1.89%  '

Set-TextValue 'D46' '2.53'
$ws.Range('E46').Value = '  -9.86%  '

Set-TextValue 'D47' '135.16'
$ws.Range('E47').Value = '  +0.00%  '

Set-TextValue 'D49' '24.78'
$ws.Range('E49').Value = '  -0.67%  '

Set-TextValue 'D50' '2.19'
$ws.Range('E50').Value = '  -1.49%  '

$ws.Range('E51').Value = '  -0.82%  '
